# Apply numeric substitutions to the multiplication worksheet.
# Each old expression is unique in the document, so a simple
# Find/Replace (ReplaceAll) for each pair is sufficient and safe.

$d = $word.ActiveDocument

$replacements = @(
    @("676×4=", "829×3="),
    @("724×4=", "402×8="),
    @("759×2=", "192×7="),
    @("721×7=", "297×8="),
    @("203×5=", "224×5="),
    @("994×8=", "124×3="),
    @("718×8=", "565×9="),
    @("492×7=", "321×8="),
    @("434×7=", "663×9="),
    @("223×2=", "513×9="),
    @("888×7=", "303×6="),
    @("498×2=", "913×6="),
    @("134×6=", "801×3="),
    @("999×6=", "261×3="),
    @("898×9=", "632×8="),
    @("421×9=", "447×5="),
    @("946×2=", "146×5="),
    @("856×3=", "247×3="),
    @("219×4=", "961×8="),
    @("325×6=", "274×5="),
    @("730×4=", "817×9="),
    @("714×5=", "223×8="),
    @("986×5=", "629×3="),
    @("489×5=", "548×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
